# Apply weekly update to Fruta/Granada sheet: rotate the variable fields
# (Fecha, Variedad, Calidad, Volumen, Precios, Unidad, Origen, Precio $/Kg, Kg/unidad)
# across rows 4,5,8,9,10,11,12,13,14,15. Rows 6 and 7 are untouched.
#
# The underlying re-sampled data forms a single 10-cycle: each row in the
# cycle below takes on the values previously held by the NEXT row in the
# cycle (wrapping from the last entry back to the first):
#   4 -> 14 -> 5 -> 9 -> 11 -> 12 -> 15 -> 8 -> 10 -> 13 -> (back to 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycle = @(4, 14, 5, 9, 11, 12, 15, 8, 10, 13)

# Columns that change: Fecha, Variedad, Calidad, Volumen, Precio min/max/prom,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot current values first (Value2 avoids COM Variant wrapper issues)
$snapshot = @{}
foreach ($r in $cycle) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

$count = $cycle.Count
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $cycle[$i]
    $sourceRow = $cycle[($i + 1) % $count]
    $sourceData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $sourceData[$c]
    }
}
